$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update unit label "m3" -> "m-3" in the two cells that used it ---
# (Updating both usages lets the now-unused "m3" shared string drop out
# and a fresh "m-3" string get appended, matching the target sharedStrings.)
$ws.Range("C10").Value = "m-3"
$ws.Range("G10").Value = "m-3"

# --- Row 13: add a mirrored mw=/value/kg-mole block in E13:G13 ---
$ws.Range("A13").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Value = "mw="

$ws.Range("B13").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Formula = "=B13"

$ws.Range("C13").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("G13").Value = "kg/mole"

# --- Row 20: simplify the I20 cross-check formula ---
$ws.Range("I20").Formula = "=I13*F9*F13*F15/(F14)"

# --- Row 22: new empty, styled cell at I22 (style copied from B22) ---
$ws.Range("B22").Copy()
$ws.Range("I22").PasteSpecial(-4122)

# --- Update the active selection to I22 ---
$null = $ws.Range("I22").Select()
